$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the existing row 476, shifting all the
# rows below (old 476-514) down to 478-516.
$ws.Rows("476:477").Insert()

# New row 476: Fecha 2021-09-22 (serial 44461), Calidad "Primera", $/caja 36 atados
$ws.Cells.Item(476, 1).Value = 6
$ws.Cells.Item(476, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(476, 3).Value = "Metropolitana"
$ws.Cells.Item(476, 4).Value = 44461
$ws.Cells.Item(476, 5).Value = 13
$ws.Cells.Item(476, 6).Value = 100112040
$ws.Cells.Item(476, 7).Value = "Cilantro"
$ws.Cells.Item(476, 8).Value = "Sin especificar"
$ws.Cells.Item(476, 9).Value = "Primera"
$ws.Cells.Item(476, 10).Value = 540
$ws.Cells.Item(476, 11).Value = 4000
$ws.Cells.Item(476, 12).Value = 4500
$ws.Cells.Item(476, 13).Value = 4213
$ws.Cells.Item(476, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(476, 15).Value = "Región Metropolitana"
$ws.Cells.Item(476, 16).Value = 117
$ws.Cells.Item(476, 17).Value = 36
$ws.Cells.Item(476, 18).Value = "Hortaliza"

# New row 477: Fecha 2021-09-22 (serial 44461), Calidad "Primera", $/docena de atados
$ws.Cells.Item(477, 1).Value = 6
$ws.Cells.Item(477, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(477, 3).Value = "Metropolitana"
$ws.Cells.Item(477, 4).Value = 44461
$ws.Cells.Item(477, 5).Value = 13
$ws.Cells.Item(477, 6).Value = 100112040
$ws.Cells.Item(477, 7).Value = "Cilantro"
$ws.Cells.Item(477, 8).Value = "Sin especificar"
$ws.Cells.Item(477, 9).Value = "Primera"
$ws.Cells.Item(477, 10).Value = 250
$ws.Cells.Item(477, 11).Value = 8000
$ws.Cells.Item(477, 12).Value = 9000
$ws.Cells.Item(477, 13).Value = 8440
$ws.Cells.Item(477, 14).Value = "`$/docena de atados"
$ws.Cells.Item(477, 15).Value = "Región Metropolitana"
$ws.Cells.Item(477, 16).Value = 2813
$ws.Cells.Item(477, 17).Value = 3
$ws.Cells.Item(477, 18).Value = "Hortaliza"
